$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1625006.6
$ws.Range("J17").Value = 1625006.6
$ws.Range("L17").Value = 4875019.800000001
$ws.Range("N17").Value = -4875355.800000001

$ws.Range("H64").Value = 50004972
$ws.Range("J64").Value = 5713.125
$ws.Range("L64").Value = 5713.125
$ws.Range("N64").Value = -6209.125

$ws.Range("H67").Value = 50004972
$ws.Range("J67").Value = 5713.125
$ws.Range("L67").Value = 5713.125
$ws.Range("N67").Value = -7429.125

$ws.Range("H88").Value = 3080.2104
$ws.Range("I88").Value = 4474.6
$ws.Range("J88").Value = 2582.2144
$ws.Range("K88").Value = 4474.6
$ws.Range("L88").Value = 2582.2144
$ws.Range("M88").Value = -4068.6
$ws.Range("N88").Value = -3394.2144

$ws.Range("H91").Value = 3080.2104
$ws.Range("I91").Value = 4474.6
$ws.Range("J91").Value = 2582.2144
$ws.Range("K91").Value = 4474.6
$ws.Range("L91").Value = 2582.2144
$ws.Range("M91").Value = -3070.6
$ws.Range("N91").Value = -5390.2144

$ws.Range("H107").Value = 2123.25
$ws.Range("I107").Value = 1997.6666
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 1997.6666
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = -77.66660000000002
$ws.Range("N107").Value = -6340

$ws.Range("H114").Value = 116665.664
$ws.Range("J114").Value = 116665.664
$ws.Range("L114").Value = 116665.664
$ws.Range("N114").Value = -125343.664

$ws.Range("H137").Value = 6024.2046
$ws.Range("I137").Value = 1356.8235
$ws.Range("K137").Value = 4070.4705
$ws.Range("M137").Value = -1520.4705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2118417
$ws.Range("I61").Value = 3150
$ws.Range("K61").Value = 3150
$ws.Range("M61").Value = -2938

$ws.Range("H74").Value = 14262.978
$ws.Range("I74").Value = 1268.0588
$ws.Range("J74").Value = 54429.09
$ws.Range("K74").Value = 1268.0588
$ws.Range("L74").Value = 54429.09
$ws.Range("M74").Value = -394.0588
$ws.Range("N74").Value = -56177.09

$ws.Range("H77").Value = 14262.978
$ws.Range("I77").Value = 1268.0588
$ws.Range("J77").Value = 54429.09
$ws.Range("K77").Value = 6340.294
$ws.Range("L77").Value = 272145.45
$ws.Range("M77").Value = -1972.294
$ws.Range("N77").Value = -280881.45

$ws.Range("H132").Value = 4278560
$ws.Range("I132").Value = 2971.3076
$ws.Range("J132").Value = 26511620
$ws.Range("K132").Value = 8913.9228
$ws.Range("L132").Value = 79534860
$ws.Range("M132").Value = -6383.9228
$ws.Range("N132").Value = -79539920

$ws.Range("H136").Value = 2118417
$ws.Range("I136").Value = 3150
$ws.Range("K136").Value = 9450
$ws.Range("M136").Value = -6900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7956796.5
$ws.Range("I20").Value = 15156753
$ws.Range("K20").Value = 15156753
$ws.Range("M20").Value = -15156506

$ws.Range("H86").Value = 90910840
$ws.Range("I86").Value = 1933.5555
$ws.Range("K86").Value = 1933.5555
$ws.Range("M86").Value = -810.5554999999999

$ws.Range("H89").Value = 90910840
$ws.Range("I89").Value = 1933.5555
$ws.Range("K89").Value = 9667.7775
$ws.Range("M89").Value = -4051.7775

$ws.Range("H94").Value = 3478.4285
$ws.Range("I94").Value = 2837.25
$ws.Range("J94").Value = 4333.3335
$ws.Range("K94").Value = 2837.25
$ws.Range("L94").Value = 4333.3335
$ws.Range("M94").Value = -2386.25
$ws.Range("N94").Value = -5235.3335

$ws.Range("H99").Value = 19986.63
$ws.Range("I99").Value = 22734.375
$ws.Range("K99").Value = 22734.375
$ws.Range("M99").Value = -21236.375

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 48888
$ws.Range("J74").Value = 48888
$ws.Range("L74").Value = 48888
$ws.Range("N74").Value = -50636

$ws.Range("H77").Value = 48888
$ws.Range("J77").Value = 48888
$ws.Range("L77").Value = 146664
$ws.Range("N77").Value = -155400

$ws.Range("H132").Value = 29413992
$ws.Range("I132").Value = 1999.2903
$ws.Range("J132").Value = 257356940
$ws.Range("K132").Value = 5997.8709
$ws.Range("L132").Value = 772070820
$ws.Range("M132").Value = -3467.8709
$ws.Range("N132").Value = -772075880

$ws.Range("H134").Value = 31256248
$ws.Range("I134").Value = 1583.1177
$ws.Range("K134").Value = 4749.3531
$ws.Range("M134").Value = -2214.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 137.66667
$ws.Range("I33").Value = 107.8
$ws.Range("K33").Value = 646.8
$ws.Range("M33").Value = -363.8

$ws.Range("H44").Value = 378
$ws.Range("I44").Value = 400
$ws.Range("J44").Value = 290
$ws.Range("K44").Value = 1200
$ws.Range("L44").Value = 870
$ws.Range("M44").Value = -802
$ws.Range("N44").Value = -1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4877.048
$ws.Range("I70").Value = 4921.727
$ws.Range("J70").Value = 4827.9
$ws.Range("K70").Value = 4921.727
$ws.Range("L70").Value = 4827.9
$ws.Range("M70").Value = -4651.727
$ws.Range("N70").Value = -5367.9

$ws.Range("H73").Value = 4877.048
$ws.Range("I73").Value = 4921.727
$ws.Range("J73").Value = 4827.9
$ws.Range("K73").Value = 4921.727
$ws.Range("L73").Value = 4827.9
$ws.Range("M73").Value = -3985.727
$ws.Range("N73").Value = -6699.9

$ws.Range("H80").Value = 12106.5
$ws.Range("I80").Value = 10110.615
$ws.Range("J80").Value = 14465.272
$ws.Range("K80").Value = 10110.615
$ws.Range("L80").Value = 14465.272
$ws.Range("M80").Value = -9112.615
$ws.Range("N80").Value = -16461.272

$ws.Range("H83").Value = 12106.5
$ws.Range("I83").Value = 10110.615
$ws.Range("J83").Value = 14465.272
$ws.Range("K83").Value = 50553.075
$ws.Range("L83").Value = 72326.36
$ws.Range("M83").Value = -45561.075
$ws.Range("N83").Value = -82310.36

$ws.Range("H102").Value = 4661342.5
$ws.Range("J102").Value = 1772.2222
$ws.Range("L102").Value = 1772.2222
$ws.Range("N102").Value = -5016.2222

$ws.Range("H122").Value = 1308281.9
$ws.Range("I122").Value = 1545787.6
$ws.Range("J122").Value = 2000.5
$ws.Range("K122").Value = 4637362.800000001
$ws.Range("L122").Value = 6001.5
$ws.Range("M122").Value = -4634912.800000001
$ws.Range("N122").Value = -10901.5

$ws.Range("H135").Value = 180000
$ws.Range("J135").Value = 180000
$ws.Range("L135").Value = 180000
$ws.Range("N135").Value = -190140

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 49241
$ws.Range("J6").Value = 49241
$ws.Range("L6").Value = 49241
$ws.Range("N6").Value = -49465

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H68").Value = 10370.4
$ws.Range("I68").Value = 12713
$ws.Range("K68").Value = 12713
$ws.Range("M68").Value = -11964

$ws.Range("H71").Value = 10370.4
$ws.Range("I71").Value = 12713
$ws.Range("K71").Value = 63565
$ws.Range("M71").Value = -59821

$ws.Range("H82").Value = 2330.6667
$ws.Range("I82").Value = 2830.8333
$ws.Range("K82").Value = 2830.8333
$ws.Range("M82").Value = -2469.8333

$ws.Range("H85").Value = 2330.6667
$ws.Range("I85").Value = 2830.8333
$ws.Range("K85").Value = 2830.8333
$ws.Range("M85").Value = -1582.8333

$ws.Range("H93").Value = 71436750
$ws.Range("I93").Value = 111117820
$ws.Range("J93").Value = 10818.6
$ws.Range("K93").Value = 111117820
$ws.Range("L93").Value = 10818.6
$ws.Range("M93").Value = -111116572
$ws.Range("N93").Value = -13314.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 478493.66
$ws.Range("I2").Value = 669635.1
$ws.Range("J2").Value = 640
$ws.Range("K2").Value = 669635.1
$ws.Range("L2").Value = 640
$ws.Range("M2").Value = -669523.1
$ws.Range("N2").Value = -864

$ws.Range("H74").Value = 46916.668
$ws.Range("J74").Value = 46916.668
$ws.Range("L74").Value = 46916.668
$ws.Range("N74").Value = -48788.668

$ws.Range("H77").Value = 46916.668
$ws.Range("J77").Value = 46916.668
$ws.Range("L77").Value = 140750.004
$ws.Range("N77").Value = -150110.004

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H96").Value = 1753.3636
$ws.Range("I96").Value = 1422.5
$ws.Range("J96").Value = 1942.4286
$ws.Range("K96").Value = 1422.5
$ws.Range("L96").Value = 1942.4286
$ws.Range("M96").Value = -49.5
$ws.Range("N96").Value = -4688.4286

$ws.Range("H122").Value = 334145.75
$ws.Range("I122").Value = 437326.5
$ws.Range("J122").Value = 5257.1875
$ws.Range("K122").Value = 1311979.5
$ws.Range("L122").Value = 15771.5625
$ws.Range("M122").Value = -1309529.5
$ws.Range("N122").Value = -20671.5625

$ws.Range("H132").Value = 245780.33
$ws.Range("J132").Value = 610549.9399999999
$ws.Range("L132").Value = 1831649.82
$ws.Range("N132").Value = -1836709.82
